$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume(1h) values as of the latest GitHub Actions run.
# For Price (column D) values that are fully numeric-looking (e.g. "37.52"), Excel's
# COM layer would otherwise coerce the assigned text into a floating point number
# (losing the literal text formatting, e.g. turning "131.60" into 131.6). To keep the
# cell a literal text string (matching the source data, which stores these as text),
# we briefly mark the cell as Text-formatted before assigning, then restore the default
# "Normal" style so no stray formatting is left behind on the cell.

$ws.Range('D2').Value = '63.640.29'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '3.089.86'
$ws.Range('E3').Value = '  -1.48%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.45%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.82%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  +2.54%  '

$ws.Range('D9').Value = '3.082.78'
$ws.Range('E9').Value = '  -1.50%  '

$ws.Range('E10').Value = '  -1.67%  '

$ws.Range('E11').Value = '  -1.03%  '

$ws.Range('E12').Value = '  -0.24%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.52'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E14').Value = '  -1.88%  '

$ws.Range('D15').Value = '3.601.50'
$ws.Range('E15').Value = '  -1.52%  '

$ws.Range('E16').Value = '  -1.55%  '

$ws.Range('E17').Value = '  -1.78%  '

$ws.Range('D18').Value = '63.581.73'
$ws.Range('E18').Value = '  -0.29%  '

$ws.Range('D19').Value = '3.085.08'
$ws.Range('E19').Value = '  -1.69%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '477.25'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.07%  '

$ws.Range('E21').Value = '  +2.07%  '

$ws.Range('E22').Value = '  -2.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.91'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.37%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.53%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.48%  '

$ws.Range('E28').Value = '  -0.16%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.98%  '

$ws.Range('E30').Value = '  -0.70%  '

$ws.Range('E31').Value = '  -0.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.30%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.80%  '

$ws.Range('E34').Value = '  -1.63%  '

$ws.Range('D35').Value = '0.0₃0851'
$ws.Range('E35').Value = '  +1.01%  '

$ws.Range('E36').Value = '  -1.48%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.53%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.97%  '

$ws.Range('E39').Value = '  -3.03%  '

$ws.Range('E40').Value = '  -0.39%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.61'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.26%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '444.14'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.96%  '

$ws.Range('E43').Value = '  -2.88%  '

$ws.Range('E44').Value = '  -2.38%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.06'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.41%  '

$ws.Range('E46').Value = '  +3.42%  '

$ws.Range('D47').Value = '2.803.01'
$ws.Range('E47').Value = '  -3.90%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.85%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.70'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.57%  '

$ws.Range('E51').Value = '  +1.30%  '
